# Supplemental-figures update: refresh the PSSM score matrix (B2:K21) with the
# recomputed values from the new run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object "object[,]" 20,10
# row 2 (F)
$data[0,0] = -15.77316496890112
$data[0,1] = 0.1657943189925281
$data[0,2] = -15.77316496890112
$data[0,3] = -15.77316496890112
$data[0,4] = -15.77316496890112
$data[0,5] = -15.77316496890112
$data[0,6] = -15.77316496890112
$data[0,7] = -15.77316496890112
$data[0,8] = -15.77316496890112
$data[0,9] = -15.77316496890112
# row 3 (W)
$data[1,0] = -15.77316496890112
$data[1,1] = -15.77316496890112
$data[1,2] = -15.77316496890112
$data[1,3] = -15.77316496890112
$data[1,4] = -15.77316496890112
$data[1,5] = -15.77316496890112
$data[1,6] = -15.77316496890112
$data[1,7] = -15.77316496890112
$data[1,8] = -15.77316496890112
$data[1,9] = -15.77316496890112
# row 4 (Y)
$data[2,0] = -15.77316496890112
$data[2,1] = 0.07772995192822135
$data[2,2] = 1.185852107222364
$data[2,3] = -15.77316496890112
$data[2,4] = 3.874351427837866
$data[2,5] = -15.77316496890112
$data[2,6] = 0.8704235084798257
$data[2,7] = -15.77316496890112
$data[2,8] = 1.356128111450722
$data[2,9] = -15.77316496890112
# row 5 (P)
$data[3,0] = -15.77316496890112
$data[3,1] = 0.6232558515837008
$data[3,2] = -15.77316496890112
$data[3,3] = -15.77316496890112
$data[3,4] = -15.77316496890112
$data[3,5] = 3.466044288193987
$data[3,6] = -15.77316496890112
$data[3,7] = -15.77316496890112
$data[3,8] = -15.77316496890112
$data[3,9] = -15.77316496890112
# row 6 (M)
$data[4,0] = -15.77316496890112
$data[4,1] = -15.77316496890112
$data[4,2] = -15.77316496890112
$data[4,3] = -15.77316496890112
$data[4,4] = -15.77316496890112
$data[4,5] = -15.77316496890112
$data[4,6] = -15.77316496890112
$data[4,7] = -15.77316496890112
$data[4,8] = -15.77316496890112
$data[4,9] = -15.77316496890112
# row 7 (I)
$data[5,0] = 3.107545543950135
$data[5,1] = -15.77316496890112
$data[5,2] = -15.77316496890112
$data[5,3] = -15.77316496890112
$data[5,4] = -15.77316496890112
$data[5,5] = -15.77316496890112
$data[5,6] = -15.77316496890112
$data[5,7] = -15.77316496890112
$data[5,8] = -15.77316496890112
$data[5,9] = -15.77316496890112
# row 8 (L)
$data[6,0] = -15.77316496890112
$data[6,1] = -15.77316496890112
$data[6,2] = -15.77316496890112
$data[6,3] = 1.662187877231675
$data[6,4] = -15.77316496890112
$data[6,5] = -15.77316496890112
$data[6,6] = -15.77316496890112
$data[6,7] = -15.77316496890112
$data[6,8] = -15.77316496890112
$data[6,9] = -15.77316496890112
# row 9 (V)
$data[7,0] = 3.508495704488655
$data[7,1] = -15.77316496890112
$data[7,2] = -15.77316496890112
$data[7,3] = -15.77316496890112
$data[7,4] = -15.77316496890112
$data[7,5] = -15.77316496890112
$data[7,6] = -15.77316496890112
$data[7,7] = -15.77316496890112
$data[7,8] = -15.77316496890112
$data[7,9] = -15.77316496890112
# row 10 (A)
$data[8,0] = -15.77316496890112
$data[8,1] = -15.77316496890112
$data[8,2] = -15.77316496890112
$data[8,3] = -15.77316496890112
$data[8,4] = -15.77316496890112
$data[8,5] = -15.77316496890112
$data[8,6] = -15.77316496890112
$data[8,7] = -15.77316496890112
$data[8,8] = -15.77316496890112
$data[8,9] = 1.328287564622785
# row 11 (G)
$data[9,0] = -15.77316496890112
$data[9,1] = -15.77316496890112
$data[9,2] = -15.77316496890112
$data[9,3] = 2.23974868228828
$data[9,4] = -15.77316496890112
$data[9,5] = 1.985481173768247
$data[9,6] = -15.77316496890112
$data[9,7] = -15.77316496890112
$data[9,8] = -15.77316496890112
$data[9,9] = 1.394066947480928
# row 12 (C)
$data[10,0] = -15.77316496890112
$data[10,1] = -15.77316496890112
$data[10,2] = -15.77316496890112
$data[10,3] = -15.77316496890112
$data[10,4] = -15.77316496890112
$data[10,5] = -15.77316496890112
$data[10,6] = -15.77316496890112
$data[10,7] = -15.77316496890112
$data[10,8] = -15.77316496890112
$data[10,9] = -15.77316496890112
# row 13 (S)
$data[11,0] = -15.77316496890112
$data[11,1] = -15.77316496890112
$data[11,2] = -15.77316496890112
$data[11,3] = 2.209210347296188
$data[11,4] = -15.77316496890112
$data[11,5] = -15.77316496890112
$data[11,6] = -15.77316496890112
$data[11,7] = -15.77316496890112
$data[11,8] = 0.9467184636097832
$data[11,9] = 1.467449368895798
# row 14 (T)
$data[12,0] = -15.77316496890112
$data[12,1] = -15.77316496890112
$data[12,2] = 1.040127727239651
$data[12,3] = -15.77316496890112
$data[12,4] = -15.77316496890112
$data[12,5] = -15.77316496890112
$data[12,6] = -15.77316496890112
$data[12,7] = -15.77316496890112
$data[12,8] = -15.77316496890112
$data[12,9] = 3.258603530185809
# row 15 (N)
$data[13,0] = -15.77316496890112
$data[13,1] = -15.77316496890112
$data[13,2] = 0.2265047866143029
$data[13,3] = -15.77316496890112
$data[13,4] = -15.77316496890112
$data[13,5] = -15.77316496890112
$data[13,6] = -15.77316496890112
$data[13,7] = -15.77316496890112
$data[13,8] = -15.77316496890112
$data[13,9] = -15.77316496890112
# row 16 (Q)
$data[14,0] = -15.77316496890112
$data[14,1] = -15.77316496890112
$data[14,2] = -15.77316496890112
$data[14,3] = -15.77316496890112
$data[14,4] = -15.77316496890112
$data[14,5] = -15.77316496890112
$data[14,6] = -15.77316496890112
$data[14,7] = -15.77316496890112
$data[14,8] = 1.159338164043632
$data[14,9] = -15.77316496890112
# row 17 (D)
$data[15,0] = -15.77316496890112
$data[15,1] = 2.035855748664325
$data[15,2] = 0.3824636374513917
$data[15,3] = -15.77316496890112
$data[15,4] = -15.77316496890112
$data[15,5] = -15.77316496890112
$data[15,6] = 0.9844464075774043
$data[15,7] = 4.321903620825934
$data[15,8] = 3.358383059398525
$data[15,9] = -15.77316496890112
# row 18 (E)
$data[16,0] = -15.77316496890112
$data[16,1] = -15.77316496890112
$data[16,2] = -15.77316496890112
$data[16,3] = -15.77316496890112
$data[16,4] = -15.77316496890112
$data[16,5] = -15.77316496890112
$data[16,6] = 1.621239108333362
$data[16,7] = -15.77316496890112
$data[16,8] = 1.595895560611966
$data[16,9] = -15.77316496890112
# row 19 (H)
$data[17,0] = -15.77316496890112
$data[17,1] = -15.77316496890112
$data[17,2] = 2.919726323789339
$data[17,3] = -15.77316496890112
$data[17,4] = -15.77316496890112
$data[17,5] = -15.77316496890112
$data[17,6] = 2.519233060576339
$data[17,7] = -15.77316496890112
$data[17,8] = -15.77316496890112
$data[17,9] = -15.77316496890112
# row 20 (K)
$data[18,0] = -15.77316496890112
$data[18,1] = 2.766346019687216
$data[18,2] = 2.492564523236869
$data[18,3] = -15.77316496890112
$data[18,4] = 2.415274167794083
$data[18,5] = -15.77316496890112
$data[18,6] = 1.775182845336371
$data[18,7] = -15.77316496890112
$data[18,8] = -15.77316496890112
$data[18,9] = 1.336063756888299
# row 21 (R)
$data[19,0] = -15.77316496890112
$data[19,1] = 2.427010460307919
$data[19,2] = -15.77316496890112
$data[19,3] = 2.904461776438625
$data[19,4] = -15.77316496890112
$data[19,5] = 2.318812689820863
$data[19,6] = 1.985871458779727
$data[19,7] = -15.77316496890112
$data[19,8] = -15.77316496890112
$data[19,9] = -15.77316496890112
$ws.Range("B2:K21").Value = $data
